# Payroll_Options.xlsx edit: replace test placeholder rows with final data
$wb = $excel.ActiveWorkbook

# --- Sheet "Rates": row 4 "Provider Test" -> "John Test",
#     add rows 5-7 for John Doe, Bob Doe, Jane Doe ---
$rates = $wb.Worksheets.Item("Rates")
$rates.Range("A4").Value = "John Test"

$rates.Range("A5").Value = "John Doe"
$rates.Range("B5").Value = 30
$rates.Range("C5").Value = 35
$rates.Range("D5").Value = 35
$rates.Range("E5").Value = 35

$rates.Range("A6").Value = "Bob Doe"
$rates.Range("B6").Value = 25
$rates.Range("C6").Value = 0
$rates.Range("D6").Value = 25
$rates.Range("E6").Value = 25

$rates.Range("A7").Value = "Jane Doe"
$rates.Range("B7").Value = 20
$rates.Range("C7").Value = 40
$rates.Range("D7").Value = 40
$rates.Range("E7").Value = 40

# --- Sheet "Billing Counselor Override": remove the second data row ---
$billing = $wb.Worksheets.Item("Billing Counselor Override")
$billing.Range("A2:C2").EntireRow.Delete()
$billing.Range("A2:C2").Select() | Out-Null

# --- Sheet "Aliases": row 2 test placeholders -> Test Value / Test Value2 ---
$aliases = $wb.Worksheets.Item("Aliases")
$aliases.Range("A2").Value = "Test Value"
$aliases.Range("B2").Value = "Test Value2"

# restore the originally active sheet/tab
$rates.Activate() | Out-Null
